# Update market price / profit data across Sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 298.625
$ws.Range("I19").Value = 274.25
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = 274.25
$ws.Range("L19").Value = 323
$ws.Range("M19").Value = -99.25
$ws.Range("N19").Value = -673

# Row 57
$ws.Range("H57").Value = 36089.668
$ws.Range("J57").Value = 43780
$ws.Range("L57").Value = 131340
$ws.Range("N57").Value = -132338

# Row 111
$ws.Range("H111").Value = 1105.75
$ws.Range("I111").Value = 779.625
$ws.Range("J111").Value = 1758
$ws.Range("K111").Value = 2338.875
$ws.Range("L111").Value = 5274
$ws.Range("M111").Value = 728.125
$ws.Range("N111").Value = -11408

# Row 113
$ws.Range("H113").Value = 3000.3076
$ws.Range("I113").Value = 2392.1428
$ws.Range("J113").Value = 3709.8333
$ws.Range("K113").Value = 2392.1428
$ws.Range("L113").Value = 3709.8333
$ws.Range("M113").Value = 861.8571999999999
$ws.Range("N113").Value = -10217.8333

# Row 125
$ws.Range("H125").Value = 1458.3
$ws.Range("I125").Value = 1100
$ws.Range("J125").Value = 1547.875
$ws.Range("K125").Value = 9900
$ws.Range("L125").Value = 13930.875
$ws.Range("M125").Value = -7440
$ws.Range("N125").Value = -18850.875

# Row 137
$ws.Range("H137").Value = 1027.8572
$ws.Range("I137").Value = 892.4375
$ws.Range("J137").Value = 1141.8948
$ws.Range("K137").Value = 2677.3125
$ws.Range("L137").Value = 3425.6844
$ws.Range("M137").Value = -127.3125
$ws.Range("N137").Value = -8525.6844

# Row 139
$ws.Range("H139").Value = 70393.336
$ws.Range("J139").Value = 70393.336
$ws.Range("L139").Value = 70393.336
$ws.Range("N139").Value = -80673.336

# Row 140
$ws.Range("H140").Value = 85022.22
$ws.Range("J140").Value = 85022.22
$ws.Range("L140").Value = 85022.22
$ws.Range("N140").Value = -95382.22

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3122.75
$ws.Range("I45").Value = 3161.5881
$ws.Range("J45").Value = 3078.7334
$ws.Range("K45").Value = 3161.5881
$ws.Range("L45").Value = 3078.7334
$ws.Range("M45").Value = -2784.5881
$ws.Range("N45").Value = -3832.7334

# Row 74
$ws.Range("H74").Value = 2952.4
$ws.Range("I74").Value = 2994.5
$ws.Range("K74").Value = 2994.5
$ws.Range("M74").Value = -2120.5

# Row 77
$ws.Range("H77").Value = 2952.4
$ws.Range("I77").Value = 2994.5
$ws.Range("K77").Value = 14972.5
$ws.Range("M77").Value = -10604.5

# Row 122
$ws.Range("H122").Value = 2581.475
$ws.Range("I122").Value = 2670.0605
$ws.Range("J122").Value = 2163.8572
$ws.Range("K122").Value = 8010.181500000001
$ws.Range("L122").Value = 6491.571599999999
$ws.Range("M122").Value = -5560.181500000001
$ws.Range("N122").Value = -11391.5716

# Row 139
$ws.Range("H139").Value = 56172.855
$ws.Range("J139").Value = 56172.855
$ws.Range("L139").Value = 56172.855
$ws.Range("N139").Value = -66452.85500000001

# Row 141
$ws.Range("H141").Value = 57322.223
$ws.Range("J141").Value = 57322.223
$ws.Range("L141").Value = 57322.223
$ws.Range("N141").Value = -67682.223

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5784.162
$ws.Range("I134").Value = 1414.625
$ws.Range("J134").Value = 13851
$ws.Range("K134").Value = 4243.875
$ws.Range("L134").Value = 41553
$ws.Range("M134").Value = -1708.875
$ws.Range("N134").Value = -46623

# Row 138
$ws.Range("H138").Value = 49610
$ws.Range("J138").Value = 49610
$ws.Range("L138").Value = 49610
$ws.Range("N138").Value = -59890

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4324.1445
$ws.Range("I31").Value = 3300.2188
$ws.Range("J31").Value = 5068.8184
$ws.Range("K31").Value = 3300.2188
$ws.Range("L31").Value = 5068.8184
$ws.Range("M31").Value = -3005.2188
$ws.Range("N31").Value = -5658.8184

# Row 34
$ws.Range("H34").Value = 4324.1445
$ws.Range("I34").Value = 3300.2188
$ws.Range("J34").Value = 5068.8184
$ws.Range("K34").Value = 3300.2188
$ws.Range("L34").Value = 5068.8184
$ws.Range("M34").Value = -3098.2188
$ws.Range("N34").Value = -5472.8184

# Row 99
$ws.Range("H99").Value = 2101.9048
$ws.Range("I99").Value = 1906.3636
$ws.Range("J99").Value = 2317
$ws.Range("K99").Value = 1906.3636
$ws.Range("L99").Value = 2317
$ws.Range("M99").Value = -408.3635999999999
$ws.Range("N99").Value = -5313

# Row 105
$ws.Range("H105").Value = 1330
$ws.Range("I105").Value = 747.1429000000001
$ws.Range("J105").Value = 2350
$ws.Range("K105").Value = 747.1429000000001
$ws.Range("L105").Value = 2350
$ws.Range("M105").Value = 999.8570999999999
$ws.Range("N105").Value = -5844

# Row 126
$ws.Range("H126").Value = 2101.9048
$ws.Range("I126").Value = 1906.3636
$ws.Range("J126").Value = 2317
$ws.Range("K126").Value = 5719.0908
$ws.Range("L126").Value = 6951
$ws.Range("M126").Value = -3249.0908
$ws.Range("N126").Value = -11891

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 167189.83
$ws.Range("I121").Value = 175.8
$ws.Range("J121").Value = 286485.56
$ws.Range("K121").Value = 527.4000000000001
$ws.Range("L121").Value = 859456.6799999999
$ws.Range("M121").Value = 782.5999999999999
$ws.Range("N121").Value = -862076.6799999999

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 18520064
$ws.Range("I126").Value = 1727.5
$ws.Range("J126").Value = 33334734
$ws.Range("K126").Value = 5182.5
$ws.Range("L126").Value = 100004202
$ws.Range("M126").Value = -2712.5
$ws.Range("N126").Value = -100009142

# Row 140
$ws.Range("H140").Value = 99769.5
$ws.Range("J140").Value = 99769.5
$ws.Range("L140").Value = 99769.5
$ws.Range("N140").Value = -110129.5

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 2785.7144
$ws.Range("J5").Value = 2785.7144
$ws.Range("L5").Value = 2785.7144
$ws.Range("N5").Value = -3011.7144

# Row 40
$ws.Range("H40").Value = 2075.8215
$ws.Range("I40").Value = 1875.4
$ws.Range("J40").Value = 2576.875
$ws.Range("K40").Value = 1875.4
$ws.Range("L40").Value = 2576.875
$ws.Range("M40").Value = -1739.4
$ws.Range("N40").Value = -2848.875

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

# Row 122
$ws.Range("H122").Value = 2555.5151
$ws.Range("I122").Value = 2320.0435
$ws.Range("J122").Value = 3097.1
$ws.Range("K122").Value = 6960.130500000001
$ws.Range("L122").Value = 9291.299999999999
$ws.Range("M122").Value = -4510.130500000001
$ws.Range("N122").Value = -14191.3

# Row 132
$ws.Range("H132").Value = 3496.6382
$ws.Range("I132").Value = 3492.2856
$ws.Range("J132").Value = 3533.2
$ws.Range("K132").Value = 10476.8568
$ws.Range("L132").Value = 10599.6
$ws.Range("M132").Value = -7946.856800000001
$ws.Range("N132").Value = -15659.6

# Row 138
$ws.Range("H138").Value = 65259.5
$ws.Range("J138").Value = 65259.5
$ws.Range("L138").Value = 65259.5
$ws.Range("N138").Value = -75539.5

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 9900
$ws.Range("J70").Value = 9900
$ws.Range("L70").Value = 9900
$ws.Range("N70").Value = -10530

# Row 73
$ws.Range("H73").Value = 9900
$ws.Range("J73").Value = 9900
$ws.Range("L73").Value = 9900
$ws.Range("N73").Value = -12084

# Row 122
$ws.Range("H122").Value = 1353.409
$ws.Range("I122").Value = 1310.2941
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3930.8823
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1480.8823
$ws.Range("N122").Value = -9400

# Row 126
$ws.Range("H126").Value = 2104.75
$ws.Range("I126").Value = 1776.25
$ws.Range("J126").Value = 2433.25
$ws.Range("K126").Value = 5328.75
$ws.Range("L126").Value = 7299.75
$ws.Range("M126").Value = -2858.75
$ws.Range("N126").Value = -12239.75

# Row 138
$ws.Range("H138").Value = 49775
$ws.Range("J138").Value = 49775
$ws.Range("L138").Value = 49775
$ws.Range("N138").Value = -60055

# Row 139
$ws.Range("H139").Value = 54683.332
$ws.Range("J139").Value = 54683.332
$ws.Range("L139").Value = 54683.332
$ws.Range("N139").Value = -64963.332

# Row 141
$ws.Range("H141").Value = 61285.715
$ws.Range("J141").Value = 63500
$ws.Range("L141").Value = 63500
$ws.Range("N141").Value = -73860

Write-Host "Market data refresh applied."
